# Actualización automática 2025-06-05 10:49:05
# Adds the "CUMPLIMIENTO MENSUAL" worksheet with per-group budget vs sales compliance data.

$wb = $excel.ActiveWorkbook

# Reference sheet used to clone existing cell styles (so we reuse the same
# style indices instead of Excel minting brand new font/xf entries).
$wsRef = $wb.Worksheets.Item("VENTA MENSUAL")

# Add the new worksheet at the end of the workbook (after the last existing sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# ---- Column widths (character units) ----
$ws.Columns.Item(1).ColumnWidth = 19.16666666666667
$ws.Columns.Item(2).ColumnWidth = 21.16666666666667
$ws.Columns.Item(3).ColumnWidth = 21.16666666666667
$ws.Columns.Item(4).ColumnWidth = 11.16666666666667
$ws.Columns.Item(5).ColumnWidth = 21.16666666666667
$ws.Columns.Item(6).ColumnWidth = 17.16666666666667

# ---- Header row (row 1) : clone formatting from an existing bold/border header ----
$wsRef.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

# ---- Data rows ----
$asesor = "CHANDI ERAZO JOSUE"

$grupos = @(
    @("240X80 PORCELANATO", 300, 0, 300, 0),
    @("FREGADEROS DE COCINA", 250.631825420901, 0, 250.631825420901, 0),
    @("GRANITO", 238.32, 0, 238.32, 0),
    @("GRIFERIAS", 106.82, 0, 106.82, 0),
    @("INODOROS", 800, 0, 800, 0),
    @("LAVABOS", 500, 0, 500, 0),
    @("LED", 300, 0, 300, 0),
    @("NO RESURTIBLES", 325.13, 0, 325.13, 0),
    @("OTROS", 0, 851.4299999999999, -851.4299999999999, 0),
    @("PANELES DECORATIVOS", 350, 0, 350, 0),
    @("PANELES PU", 230, 0, 230, 0),
    @("PANELES PVC", 483, 0, 483, 0),
    @("PORCELANATO", 7774.1, 0, 7774.1, 0),
    @("PUERTAS DE SEGURIDAD", 342, 0, 342, 0),
    @("SAL SOLUBLE", 1500, 0, 1500, 0)
)

$row = 2
foreach ($g in $grupos) {
    $ws.Cells.Item($row, 1).Value = $asesor
    $ws.Cells.Item($row, 2).Value = $g[0]

    # Money columns (PRESUPUESTO / VENTA / POR CUMPLIR) - reuse the "$#,##0.00" style.
    $wsRef.Range("C2").Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
    $ws.Cells.Item($row, 3).Value = $g[1]

    $wsRef.Range("C2").Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)
    $ws.Cells.Item($row, 4).Value = $g[2]

    $wsRef.Range("C2").Copy()
    $ws.Cells.Item($row, 5).PasteSpecial(-4122)
    $ws.Cells.Item($row, 5).Value = $g[3]

    # Percentage column (CUMPLIMIENTO).
    $ws.Cells.Item($row, 6).NumberFormat = "0.00%"
    $ws.Cells.Item($row, 6).Value = $g[4]

    $row = $row + 1
}

# ---- Totals row (row 17) ----
$ws.Range("B17").Value = "TOTAL"
$ws.Range("B17").HorizontalAlignment = -4152

$wsRef.Range("C2").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 13500.0018254209

$wsRef.Range("C2").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 851.4299999999999

$wsRef.Range("C2").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = 12648.5718254209

$ws.Range("F17").NumberFormat = "0.00%"
$ws.Range("F17").Value = 0.06306888036094425

$excel.CutCopyMode = 0

# Restore the originally active sheet/tab.
$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select()

Write-Host "CUMPLIMIENTO MENSUAL sheet added"
